# Sync automático del tracker (cada 3h)
# Appends the two newest rows synced from the results feed to the bottom
# of the tracker sheet:
#   event_id | fecha | jugador_A | jugador_B | pronostico | cuota | resultado | profit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 139; EventId = "14678166"; Fecha = "2025-09-15"; JugadorA = "Lois Boisson"; JugadorB = "Yeonwoo Ku";   Pronostico = "Gana Lois Boisson"; Cuota = 1.36 },
    @{ Row = 140; EventId = "14679464"; Fecha = "2025-09-14"; JugadorA = "Nico Hipfl";   JugadorB = "Marvin Möller"; Pronostico = "Gana Nico Hipfl";   Cuota = 4.33 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # event_id is synced from the feed as a text code (not a calculated
    # number), so the cell is pre-formatted as Text before the value lands.
    $eventCell = $ws.Cells.Item($row, 1)
    $eventCell.NumberFormat = "@"
    $eventCell.Value = $r.EventId

    # fecha is stored as plain "yyyy-mm-dd" text, same as every other row.
    $fechaCell = $ws.Cells.Item($row, 2)
    $fechaCell.NumberFormat = "@"
    $fechaCell.Value = $r.Fecha

    $ws.Cells.Item($row, 3).Value = $r.JugadorA
    $ws.Cells.Item($row, 4).Value = $r.JugadorB
    $ws.Cells.Item($row, 5).Value = $r.Pronostico
    $ws.Cells.Item($row, 6).Value = $r.Cuota

    # resultado / profit are still unknown for these just-synced matches —
    # leave them as the same blank text placeholders used by every other
    # pending row (pre-formatted as Text so the cell stays materialized
    # instead of disappearing as a truly empty cell).
    $resultCell = $ws.Cells.Item($row, 7)
    $resultCell.NumberFormat = "@"
    $resultCell.Value = ""

    $profitCell = $ws.Cells.Item($row, 8)
    $profitCell.NumberFormat = "@"
    $profitCell.Value = ""
}
